# Updates cryptos list data (B/C/D/E columns) per the Feb 23 2023 scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.462.23"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.65"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9977"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.23"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9974"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3948"
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.04"
$ws.Range("E9").Value = "  +5.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.407"
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9975"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08587"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.55"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.338"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001337"
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.859"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.664.16"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.64"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06961"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.61"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.012"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9965"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.72"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.462.22"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.437"
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.038"
$ws.Range("E26").Value = "  +10.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.55"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.93"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.93"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.347"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.111"
$ws.Range("E31").Value = "  -5.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.537"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.840.66"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.074"
$ws.Range("E34").Value = "  +8.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08257"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.25"
$ws.Range("E36").Value = "  +12.02%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02992"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.830"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2759"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09271"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7764"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.87"
$ws.Range("E42").Value = "  +5.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.444"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.58"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7138"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.533"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9968"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.40"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.447"
$ws.Range("E51").Value = "  +12.15%  "
